# Add new workout rows (380-386) to the Kilimanjaro Weekly Scoreboard sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Preserve the short-date number format already used throughout column B
# by copying the formatting (only) from the last existing data row.
$ws.Range("B379").Copy()
$ws.Range("B380:B386").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# New data rows: Participant, Date (serial), Workout Type, Total Duration,
# Total Distance, Total Elevation, Zone1..Zone5, Workout Level, Week
$newRows = @(
    @("Steven",    45518, "Workout", 26, 0,    0,   22, 4,  0,  0, 0, "Brave Leopard",      10),
    @("Eric",      45518, "Workout", 94, 0,    0,   40, 46, 8,  0, 0, "Sauntering Hippo",   10),
    @("Steven",    45518, "Walk",    35, 1.79, 108, 35, 0,  0,  0, 0, "Brave Leopard",      10),
    @("Matt",      45518, "Walk",    34, 1.44, 108, 34, 0,  0,  0, 0, "Agile Antelope",     10),
    @("Jeremiah",  45519, "Run",     20, 2.07, 151, 0,  13, 2,  0, 0, "Sauntering Hippo",   10),
    @("Steven",    45519, "Walk",    23, 1.15, 138, 23, 0,  0,  0, 0, "Brave Leopard",      10),
    @("Steven",    45519, "Run",     29, 2.68, 154, 2,  2,  24, 1, 0, "Brave Leopard",      10)
)

$r = 380
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value  = $row[0]    # A Participant
    $ws.Cells.Item($r, 2).Value  = $row[1]    # B Date
    $ws.Cells.Item($r, 3).Value  = $row[2]    # C Workout Type
    $ws.Cells.Item($r, 4).Value  = $row[3]    # D Total Duration
    $ws.Cells.Item($r, 5).Value  = $row[4]    # E Total Distance
    $ws.Cells.Item($r, 6).Value  = $row[5]    # F Total Elevation
    $ws.Cells.Item($r, 7).Value  = $row[6]    # G Zone 1
    $ws.Cells.Item($r, 8).Value  = $row[7]    # H Zone 2
    $ws.Cells.Item($r, 9).Value  = $row[8]    # I Zone 3
    $ws.Cells.Item($r, 10).Value = $row[9]    # J Zone 4
    $ws.Cells.Item($r, 11).Value = $row[10]   # K Zone 5
    $ws.Cells.Item($r, 12).Value = $row[11]   # L Workout Level
    $ws.Cells.Item($r, 13).Value = $row[12]   # M Week
    $r++
}

# Scroll the frozen view down and select the first empty cell beneath the
# newly added data, matching where the author left off after pasting.
$ws.Range("D387").Select()
$excel.ActiveWindow.ScrollRow = 362
